$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column so numeric-looking strings
# (e.g. "2.00", "7.04") are not auto-converted to numbers, matching the
# original inlineStr text cells in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '62.739.83'
$ws.Range("E2").Value = '  +2.31%  '

# Row 3
$ws.Range("D3").Value = '2.939.47'
$ws.Range("E3").Value = '  +0.23%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = '591.51'
$ws.Range("E5").Value = '  -0.59%  '

# Row 6
$ws.Range("D6").Value = '146.89'
$ws.Range("E6").Value = '  +2.34%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").Value = '2.938.30'
$ws.Range("E8").Value = '  +0.24%  '

# Row 9
$ws.Range("E9").Value = '  +0.75%  '

# Row 10
$ws.Range("D10").Value = '7.04'
$ws.Range("E10").Value = '  +1.55%  '

# Row 11
$ws.Range("E11").Value = '  +5.01%  '

# Row 12
$ws.Range("D12").Value = '0.436'
$ws.Range("E12").Value = '  -0.16%  '

# Row 13
$ws.Range("D13").Value = '0.0000232'
$ws.Range("E13").Value = '  +3.97%  '

# Row 14
$ws.Range("D14").Value = '32.38'
$ws.Range("E14").Value = '  -2.61%  '

# Row 15
$ws.Range("E15").Value = '  -1.02%  '

# Row 16
$ws.Range("D16").Value = '3.431.94'
$ws.Range("E16").Value = '  +0.42%  '

# Row 17
$ws.Range("D17").Value = '62.758.84'
$ws.Range("E17").Value = '  +2.38%  '

# Row 18
$ws.Range("E18").Value = '  +0.25%  '

# Row 19
$ws.Range("D19").Value = '2.942.42'
$ws.Range("E19").Value = '  +0.24%  '

# Row 20
$ws.Range("D20").Value = '437.64'
$ws.Range("E20").Value = '  +0.91%  '

# Row 21
$ws.Range("D21").Value = '13.37'
$ws.Range("E21").Value = '  -1.50%  '

# Row 22
$ws.Range("D22").Value = '0.663'
$ws.Range("E22").Value = '  -1.42%  '

# Row 23
$ws.Range("D23").Value = '6.98'
$ws.Range("E23").Value = '  -1.10%  '

# Row 24
$ws.Range("D24").Value = '11.14'
$ws.Range("E24").Value = '  +2.68%  '

# Row 25
$ws.Range("D25").Value = '80.68'
$ws.Range("E25").Value = '  -0.96%  '

# Row 26
$ws.Range("D26").Value = '11.77'
$ws.Range("E26").Value = '  +0.48%  '

# Row 27
$ws.Range("D27").Value = '2.11'
$ws.Range("E27").Value = '  -2.59%  '

# Row 28
$ws.Range("E28").Value = '  -0.05%  '

# Row 29
$ws.Range("E29").Value = '  +0.97%  '

# Row 30
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  +5.38%  '

# Row 31
$ws.Range("E31").Value = '  +0.23%  '

# Row 32
$ws.Range("D32").Value = '0.0₃0979'
$ws.Range("E32").Value = '  +12.02%  '

# Row 33
$ws.Range("D33").Value = '26.27'
$ws.Range("E33").Value = '  -1.46%  '

# Row 34
$ws.Range("E34").Value = '  -0.70%  '

# Row 35
$ws.Range("E35").Value = '  +0.06%  '

# Row 36
$ws.Range("D36").Value = '0.989'
$ws.Range("E36").Value = '  -2.38%  '

# Row 37
$ws.Range("D37").Value = '5.59'
$ws.Range("E37").Value = '  -0.49%  '

# Row 38
$ws.Range("D38").Value = '2.99'
$ws.Range("E38").Value = '  +0.73%  '

# Row 39
$ws.Range("D39").Value = '49.57'
$ws.Range("E39").Value = '  -0.38%  '

# Row 40
$ws.Range("D40").Value = '2.00'
$ws.Range("E40").Value = '  +0.95%  '

# Row 41
$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").Value = '8.42'
$ws.Range("E41").Value = '  -0.88%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").Value = '  -4.06%  '

# Row 43
$ws.Range("D43").Value = '0.277'
$ws.Range("E43").Value = '  -0.35%  '

# Row 44
$ws.Range("D44").Value = '39.12'
$ws.Range("E44").Value = '  -7.05%  '

# Row 45
$ws.Range("D45").Value = '2.698.97'
$ws.Range("E45").Value = '  -0.20%  '

# Row 46
$ws.Range("D46").Value = '135.08'
$ws.Range("E46").Value = '  +1.27%  '

# Row 47
$ws.Range("D47").Value = '0.0336'
$ws.Range("E47").Value = '  -2.12%  '

# Row 48
$ws.Range("D48").Value = '355.87'
$ws.Range("E48").Value = '  -2.23%  '

# Row 50
$ws.Range("E50").Value = '  -0.77%  '

# Row 51
$ws.Range("D51").Value = '22.56'
$ws.Range("E51").Value = '  -4.22%  '
